$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 1853.3334
$ws.Range("I116").Value = 1836.3636
$ws.Range("K116").Value = 1836.3636
$ws.Range("M116").Value = 1605.6364

$ws.Range("H132").Value = 5717818
$ws.Range("I132").Value = 6806377.5
$ws.Range("J132").Value = 2879.25
$ws.Range("K132").Value = 20419132.5
$ws.Range("L132").Value = 8637.75
$ws.Range("M132").Value = -20416602.5
$ws.Range("N132").Value = -13697.75

$ws.Range("H137").Value = 1324.4103
$ws.Range("I137").Value = 998.25
$ws.Range("J137").Value = 2154.6365
$ws.Range("K137").Value = 2994.75
$ws.Range("L137").Value = 6463.9095
$ws.Range("M137").Value = -444.75
$ws.Range("N137").Value = -11563.9095

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1451
$ws.Range("I2").Value = 1133.5
$ws.Range("J2").Value = 1689.125
$ws.Range("K2").Value = 1133.5
$ws.Range("L2").Value = 1689.125
$ws.Range("M2").Value = -1020.5
$ws.Range("N2").Value = -1915.125

$ws.Range("H45").Value = 1255.8422
$ws.Range("I45").Value = 1229.2727
$ws.Range("J45").Value = 1292.375
$ws.Range("K45").Value = 1229.2727
$ws.Range("L45").Value = 1292.375
$ws.Range("M45").Value = -852.2727
$ws.Range("N45").Value = -2046.375

$ws.Range("H61").Value = 1433.2041
$ws.Range("I61").Value = 1047.6316
$ws.Range("J61").Value = 2765.182
$ws.Range("K61").Value = 1047.6316
$ws.Range("L61").Value = 2765.182
$ws.Range("M61").Value = -835.6315999999999
$ws.Range("N61").Value = -3189.182

$ws.Range("H74").Value = 992.04346
$ws.Range("I74").Value = 835.15
$ws.Range("J74").Value = 2038
$ws.Range("K74").Value = 835.15
$ws.Range("L74").Value = 2038
$ws.Range("M74").Value = 38.85000000000002
$ws.Range("N74").Value = -3786

$ws.Range("H77").Value = 992.04346
$ws.Range("I77").Value = 835.15
$ws.Range("J77").Value = 2038
$ws.Range("K77").Value = 4175.75
$ws.Range("L77").Value = 10190
$ws.Range("M77").Value = 192.25
$ws.Range("N77").Value = -18926

$ws.Range("H116").Value = 1451
$ws.Range("I116").Value = 1133.5
$ws.Range("J116").Value = 1689.125
$ws.Range("K116").Value = 1133.5
$ws.Range("L116").Value = 1689.125
$ws.Range("M116").Value = 1160.5
$ws.Range("N116").Value = -6277.125

$ws.Range("H136").Value = 1433.2041
$ws.Range("I136").Value = 1047.6316
$ws.Range("J136").Value = 2765.182
$ws.Range("K136").Value = 3142.8948
$ws.Range("L136").Value = 8295.545999999998
$ws.Range("M136").Value = -592.8948
$ws.Range("N136").Value = -13395.546

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1451
$ws.Range("I3").Value = 1133.5
$ws.Range("J3").Value = 1689.125
$ws.Range("K3").Value = 1133.5
$ws.Range("L3").Value = 1689.125
$ws.Range("M3").Value = -1019.5
$ws.Range("N3").Value = -1917.125

$ws.Range("H99").Value = 1199.8
$ws.Range("I99").Value = 1249.75
$ws.Range("J99").Value = 1000
$ws.Range("K99").Value = 1249.75
$ws.Range("L99").Value = 1000
$ws.Range("M99").Value = 248.25
$ws.Range("N99").Value = -3996

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 894.16
$ws.Range("I58").Value = 745.8261
$ws.Range("J58").Value = 2600
$ws.Range("K58").Value = 745.8261
$ws.Range("L58").Value = 2600
$ws.Range("M58").Value = -542.8261
$ws.Range("N58").Value = -3006

$ws.Range("H99").Value = 32628.666
$ws.Range("I99").Value = 73445.36
$ws.Range("J99").Value = 2553.2104
$ws.Range("K99").Value = 73445.36
$ws.Range("L99").Value = 2553.2104
$ws.Range("M99").Value = -71947.36
$ws.Range("N99").Value = -5549.2104

$ws.Range("H126").Value = 32628.666
$ws.Range("I126").Value = 73445.36
$ws.Range("J126").Value = 2553.2104
$ws.Range("K126").Value = 220336.08
$ws.Range("L126").Value = 7659.6312
$ws.Range("M126").Value = -217866.08
$ws.Range("N126").Value = -12599.6312

$ws.Range("H134").Value = 1296.3684
$ws.Range("I134").Value = 1279.4286
$ws.Range("J134").Value = 1343.8
$ws.Range("K134").Value = 3838.2858
$ws.Range("L134").Value = 4031.4
$ws.Range("M134").Value = -1303.2858
$ws.Range("N134").Value = -9101.4

$ws.Range("H136").Value = 894.16
$ws.Range("I136").Value = 745.8261
$ws.Range("J136").Value = 2600
$ws.Range("K136").Value = 2237.4783
$ws.Range("L136").Value = 7800
$ws.Range("M136").Value = 312.5217000000002
$ws.Range("N136").Value = -12900

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H21").Value = 156.14285
$ws.Range("I21").Value = 156.14285
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 468.42855
$ws.Range("L21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("N21").Value = -295.42855

$ws.Range("H22").Value = 4999
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()

$ws.Range("H27").Value = 4999
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()

$ws.Range("H116").Value = 5777
$ws.Range("I116").Value = 6957.3125
$ws.Range("J116").Value = 2000
$ws.Range("K116").Value = 20871.9375
$ws.Range("L116").Value = 6000
$ws.Range("M116").Value = -17429.9375
$ws.Range("N116").Value = -12884

$ws.Range("H122").Value = 1439.4667
$ws.Range("I122").Value = 1538
$ws.Range("J122").Value = 1326.8572
$ws.Range("K122").Value = 13842
$ws.Range("L122").Value = 11941.7148
$ws.Range("M122").Value = -11392
$ws.Range("N122").Value = -16841.7148

$ws.Range("H131").Value = 2546.125
$ws.Range("J131").Value = 938.3714
$ws.Range("L131").Value = 2815.1142
$ws.Range("N131").Value = -12895.1142

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 31090.908
$ws.Range("I3").Value = 34666.668
$ws.Range("J3").Value = 15000
$ws.Range("K3").Value = 34666.668
$ws.Range("L3").Value = 15000
$ws.Range("M3").Value = -34550.668
$ws.Range("N3").Value = -15232

$ws.Range("H9").Value = 215
$ws.Range("I9").Value = 200
$ws.Range("J9").Value = 230
$ws.Range("K9").Value = 200
$ws.Range("L9").Value = 230
$ws.Range("M9").Value = -30
$ws.Range("N9").Value = -570

$ws.Range("H11").Value = 2335874
$ws.Range("I11").Value = 2145419
$ws.Range("J11").Value = 3002466.8
$ws.Range("K11").Value = 2145419
$ws.Range("L11").Value = 3002466.8
$ws.Range("M11").Value = -2145280
$ws.Range("N11").Value = -3002744.8

$ws.Range("H13").Value = 4521
$ws.Range("J13").Value = 9750
$ws.Range("L13").Value = 9750
$ws.Range("N13").Value = -10028

$ws.Range("H70").Value = 12754063
$ws.Range("I70").Value = 13425119
$ws.Range("K70").Value = 13425119
$ws.Range("M70").Value = -13424849

$ws.Range("H73").Value = 12754063
$ws.Range("I73").Value = 13425119
$ws.Range("K73").Value = 13425119
$ws.Range("M73").Value = -13424183

$ws.Range("H95").Value = 10000
$ws.Range("J95").Value = 10000
$ws.Range("L95").Value = 10000
$ws.Range("N95").Value = -15492

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 4358.0264
$ws.Range("I136").Value = 5341.4287
$ws.Range("J136").Value = 1604.5
$ws.Range("K136").Value = 16024.2861
$ws.Range("L136").Value = 4813.5
$ws.Range("M136").Value = -13474.2861
$ws.Range("N136").Value = -9913.5
